$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.335.86'
$ws.Range('E2').Value = '  +1.24%  '
$ws.Range('D3').Value = '2.010.49'
$ws.Range('E3').Value = '  +5.05%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '324.79'
$ws.Range('E5').Value = '  +1.53%  '
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').Value = '0.5128'
$ws.Range('E7').Value = '  +1.74%  '
$ws.Range('D8').Value = '0.4263'
$ws.Range('E8').Value = '  +5.48%  '
$ws.Range('D9').Value = '0.08703'
$ws.Range('E9').Value = '  +5.10%  '
$ws.Range('E10').Value = '  +3.01%  '
$ws.Range('D11').Value = '43.18'
$ws.Range('D12').Value = '24.76'
$ws.Range('E12').Value = '  +2.86%  '
$ws.Range('D13').Value = '2.009.08'
$ws.Range('E13').Value = '  +4.58%  '
$ws.Range('D14').Value = '6.573'
$ws.Range('E14').Value = '  +2.88%  '
$ws.Range('D15').Value = '7.474'
$ws.Range('E15').Value = '  +3.56%  '
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').Value = '94.38'
$ws.Range('E17').Value = '  +2.39%  '
$ws.Range('E18').Value = '  +1.71%  '
$ws.Range('D19').Value = '0.06540'
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D20').Value = '18.87'
$ws.Range('E20').Value = '  +3.91%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '6.200'
$ws.Range('E22').Value = '  +4.51%  '
$ws.Range('D23').Value = '30.396.36'
$ws.Range('E23').Value = '  +1.29%  '
$ws.Range('D24').Value = '11.81'
$ws.Range('E24').Value = '  +4.53%  '
$ws.Range('D25').Value = '2.260'
$ws.Range('E25').Value = '  +2.95%  '
$ws.Range('D26').Value = '2.245.79'
$ws.Range('E26').Value = '  +4.96%  '
$ws.Range('D27').Value = '22.43'
$ws.Range('E27').Value = '  +1.26%  '
$ws.Range('D28').Value = '162.27'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').Value = '2.423'
$ws.Range('E29').Value = '  +5.47%  '
$ws.Range('D30').Value = '131.09'
$ws.Range('E30').Value = '  +1.62%  '
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('E32').Value = '  +1.64%  '
$ws.Range('D33').Value = '6.079'
$ws.Range('E33').Value = '  +2.40%  '
$ws.Range('D34').Value = '3.826'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').Value = '1.370'
$ws.Range('E35').Value = '  +14.60%  '
$ws.Range('D36').Value = '0.02524'
$ws.Range('E36').Value = '  +3.29%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '0.06687'
$ws.Range('E37').Value = '  +4.56%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '5.469'
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').Value = '12.37'
$ws.Range('E39').Value = '  +8.67%  '
$ws.Range('D40').Value = '9.115'
$ws.Range('E40').Value = '  +4.47%  '
$ws.Range('D41').Value = '0.2193'
$ws.Range('E41').Value = '  +1.93%  '
$ws.Range('D42').Value = '0.6639'
$ws.Range('E42').Value = '  +2.73%  '
$ws.Range('D43').Value = '1.240'
$ws.Range('E43').Value = '  +2.33%  '
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').Value = '13.67'
$ws.Range('E45').Value = '  +3.13%  '
$ws.Range('D46').Value = '0.6169'
$ws.Range('E46').Value = '  +2.16%  '
$ws.Range('D47').Value = '2.180'
$ws.Range('E47').Value = '  -1.56%  '
$ws.Range('D48').Value = '3.660'
$ws.Range('E48').Value = '  +0.70%  '
$ws.Range('E49').Value = '  +4.69%  '
$ws.Range('D50').Value = '124.25'
$ws.Range('E50').Value = '  +1.67%  '
$ws.Range('D51').Value = '80.63'
$ws.Range('E51').Value = '  +2.16%  '
